$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("B2","B3","B4","C2","C3","C4","D2","D3","D4")
foreach ($c in $cells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("B2").Value = "0.25"
$ws.Range("B3").Value = "0.12"
$ws.Range("B4").Value = "-0.1*"
$ws.Range("C2").Value = "0.4***"
$ws.Range("C3").Value = "-0.57***"
$ws.Range("C4").Value = "0.02"
$ws.Range("D2").Value = "-11.46***"
$ws.Range("D3").Value = "10.16"
$ws.Range("D4").Value = "1.48*"

foreach ($c in $cells) {
    $ws.Range($c).Style = "Normal"
}
